# Update the "Förändrad" (changed) date column (C) for rows 2 through 28
# from 2024-08-24 (serial 45528) to 2024-08-25 (serial 45529).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45528) {
        $cell.Value2 = 45529
    }
}
